# Update the dSF (column F) values on Sheet1 with the newly computed
# data-pull results. Rows 14 and 18 are intentionally left unchanged
# (their dSF value already matched dS0 = 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -3
    3  = 1
    4  = 7
    5  = 3
    6  = -4
    7  = 4
    8  = -1
    9  = 1
    10 = 3
    11 = 5
    12 = 2
    13 = 8
    15 = -2
    16 = 4
    17 = -1
    19 = 1
    20 = 4
    21 = -4
    22 = -4
    23 = 2
    24 = -6
    25 = 2
    26 = 4
    27 = 5
    28 = -1
    29 = 2
    30 = -5
    31 = 4
    32 = 4
    33 = 2
    34 = 9
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
